# Adds the new "2024-2025" sheet with ex-parte order data for
# administracion_de_tribunales, matching the other three sheets' layout.

$wb = $excel.ActiveWorkbook

# --- Add the new sheet after the last existing one ("2022-2023") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2024-2025"

# --- Headers (row 1) ---
$headers = @(
    "Región",
    "Cantidad de órdenes ex parte",
    "Delito de agresión sexual en órdenes ex parte emitidas",
    "Delitos de acoso sexual en órdenes ex parte emitidas",
    "Delitos de actos lascivos en órdenes ex parte emitidas",
    "Delito de incesto en órdenes ex parte emitidas"
)
for ($c = 0; $c -lt $headers.Count; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# --- Region rows (2-14) + Total row (15) ---
# Each entry: Región, Cantidad, Agresión sexual, Acoso sexual, Actos lascivos, Incesto
$rows = @(
    ,@("Aguadilla", 7, 4, 2, 4, $null)
    ,@("Aibonito", $null, $null, $null, $null, $null)
    ,@("Arecibo", 9, $null, $null, 1, 1)
    ,@("Bayamón", 10, 4, 3, 3, $null)
    ,@("Caguas", 12, 1, 1, 1, $null)
    ,@("Carolina", 4, 1, 3, 3, 1)
    ,@("Fajardo", 1, 1, $null, $null, $null)
    ,@("Guayama", 1, 1, $null, $null, $null)
    ,@("Humacao", 1, 1, $null, $null, $null)
    ,@("Mayagüez", 4, 1, 3, 1, $null)
    ,@("Ponce", 12, 7, 1, 8, 1)
    ,@("San Juan", 3, 2, $null, $null, $null)
    ,@("Utuado", 5, 2, 3, 3, $null)
    ,@("Total", 71, 25, 16, 24, 3)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($excelRow, $c + 1).Value = $val
        }
    }
}

# --- Formatting: black font text across the populated range (A1:F15) ---
$dataRange = $ws.Range("A1:F15")
$dataRange.Font.Color = 0
$dataRange.Font.Name = "Calibri"

# --- Make the new sheet the active tab, matching the saved workbook state ---
$ws.Activate()
